$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.790.30'
$ws.Range('E2').Value = '  +8.41%  '

$ws.Range('D3').Value = '3.499.51'
$ws.Range('E3').Value = '  +11.88%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = "'188.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +12.93%  '

$ws.Range('D6').Value = "'549.67"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.95%  '

$ws.Range('D7').Value = '3.491.06'
$ws.Range('E7').Value = '  +11.85%  '

$ws.Range('D8').Value = "'0.607"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.91%  '

$ws.Range('E9').Value = '  -0.12%  '

$ws.Range('D10').Value = "'0.633"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.52%  '

$ws.Range('D11').Value = "'0.151"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +19.66%  '

$ws.Range('D12').Value = "'55.18"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.02%  '

$ws.Range('E13').Value = '  +9.41%  '

$ws.Range('D14').Value = "'9.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.69%  '

$ws.Range('D15').Value = '4.057.85'
$ws.Range('E15').Value = '  +11.83%  '

$ws.Range('D16').Value = '3.493.45'
$ws.Range('E16').Value = '  +11.71%  '

$ws.Range('E17').Value = '  +8.24%  '

$ws.Range('D18').Value = '66.796.19'
$ws.Range('E18').Value = '  +8.63%  '

$ws.Range('D19').Value = "'18.24"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.41%  '

$ws.Range('D20').Value = "'11.77"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.66%  '

$ws.Range('D21').Value = "'0.992"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.40%  '

$ws.Range('D22').Value = "'410.79"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +15.70%  '

$ws.Range('D23').Value = "'85.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.18%  '

$ws.Range('D24').Value = "'3.91"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.38%  '

$ws.Range('D25').Value = "'4.24"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.32%  '

$ws.Range('D26').Value = "'11.11"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.76%  '

$ws.Range('D27').Value = "'2.93"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.14%  '

$ws.Range('D28').Value = "'6.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.75%  '

$ws.Range('D29').Value = "'11.87"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.65%  '

$ws.Range('E30').Value = '  +11.82%  '

$ws.Range('E31').Value = '  +9.92%  '

$ws.Range('D32').Value = "'653.14"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.08%  '

$ws.Range('D33').Value = "'6.70"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.35%  '

$ws.Range('D34').Value = "'11.70"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.25%  '

$ws.Range('E35').Value = '  +10.38%  '

$ws.Range('D36').Value = "'59.56"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.78%  '

$ws.Range('D37').Value = "'38.68"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.15%  '

$ws.Range('D38').Value = '0.0₃0813'
$ws.Range('E38').Value = '  +20.43%  '

$ws.Range('E39').Value = '  -0.09%  '

$ws.Range('E40').Value = '  +8.22%  '

$ws.Range('E41').Value = '  +15.18%  '

$ws.Range('D42').Value = "'3.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +20.01%  '

$ws.Range('D43').Value = "'0.998"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.04%  '

$ws.Range('D44').Value = '3.013.84'
$ws.Range('E44').Value = '  +8.32%  '

$ws.Range('D45').Value = "'2.91"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +18.16%  '

$ws.Range('E46').Value = '  +9.43%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = "'0.0415"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.80%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = "'3.24"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +13.08%  '

$ws.Range('D49').Value = "'2.70"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.48%  '

$ws.Range('D50').Value = "'8.90"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +21.46%  '

$ws.Range('E51').Value = '  +9.14%  '
